# modifs gantt pour test et refusion
# Copy the first 5 task labels (rows 11-15, col A) down to rows 49-53
# of the Gantt sheet, then leave the selection/scroll on the newly
# added block (as seen in the edited workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A49").Value = $ws.Range("A11").Value()
$ws.Range("A50").Value = $ws.Range("A12").Value()
$ws.Range("A51").Value = $ws.Range("A13").Value()
$ws.Range("A52").Value = $ws.Range("A14").Value()
$ws.Range("A53").Value = $ws.Range("A15").Value()

# Reflect the view state captured in the edited file: the sheet is
# scrolled down and the newly populated range is selected.
$ws.Activate()
$ws.Range("A37").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A49:A53").Select()
